$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2242.8572
$ws.Range("I19").Value = 1300
$ws.Range("J19").Value = 2400
$ws.Range("K19").Value = 1300
$ws.Range("L19").Value = 2400
$ws.Range("M19").Value = -1125
$ws.Range("N19").Value = -2750
$ws.Range("H33").Value = 616.08
$ws.Range("I33").Value = 533.9048
$ws.Range("J33").Value = 1047.5
$ws.Range("K33").Value = 533.9048
$ws.Range("L33").Value = 1047.5
$ws.Range("M33").Value = -304.9048
$ws.Range("N33").Value = -1505.5
$ws.Range("H116").Value = 46883
$ws.Range("I116").Value = 58300.633
$ws.Range("J116").Value = 3496
$ws.Range("K116").Value = 58300.633
$ws.Range("L116").Value = 3496
$ws.Range("M116").Value = -54858.633
$ws.Range("N116").Value = -10380
$ws.Range("H124").Value = 33316.668
$ws.Range("J124").Value = 33316.668
$ws.Range("L124").Value = 33316.668
$ws.Range("N124").Value = -43136.668
$ws.Range("H126").Value = 32430
$ws.Range("J126").Value = 32430
$ws.Range("L126").Value = 32430
$ws.Range("N126").Value = -42310
$ws.Range("H137").Value = 6805.5356
$ws.Range("I137").Value = 10234.4375
$ws.Range("J137").Value = 2233.6667
$ws.Range("K137").Value = 30703.3125
$ws.Range("L137").Value = 6701.000100000001
$ws.Range("M137").Value = -28153.3125
$ws.Range("N137").Value = -11801.0001
$ws.Range("H140").Value = 34378
$ws.Range("J140").Value = 34378
$ws.Range("L140").Value = 34378
$ws.Range("N140").Value = -44738

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 276525.97
$ws.Range("I61").Value = 229051.16
$ws.Range("J61").Value = 348556.7
$ws.Range("K61").Value = 229051.16
$ws.Range("L61").Value = 348556.7
$ws.Range("M61").Value = -228839.16
$ws.Range("N61").Value = -348980.7
$ws.Range("H74").Value = 270315.84
$ws.Range("I74").Value = 401185.4
$ws.Range("K74").Value = 401185.4
$ws.Range("M74").Value = -400311.4
$ws.Range("H77").Value = 270315.84
$ws.Range("I77").Value = 401185.4
$ws.Range("K77").Value = 2005927
$ws.Range("M77").Value = -2001559
$ws.Range("H132").Value = 25488.445
$ws.Range("I132").Value = 43424.2
$ws.Range("J132").Value = 3068.75
$ws.Range("K132").Value = 130272.6
$ws.Range("L132").Value = 9206.25
$ws.Range("M132").Value = -127742.6
$ws.Range("N132").Value = -14266.25
$ws.Range("H136").Value = 276525.97
$ws.Range("I136").Value = 229051.16
$ws.Range("J136").Value = 348556.7
$ws.Range("K136").Value = 687153.48
$ws.Range("L136").Value = 1045670.1
$ws.Range("M136").Value = -684603.48
$ws.Range("N136").Value = -1050770.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 40000
$ws.Range("I99").Value = 50000
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 50000
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -48502
$ws.Range("N99").Value = -32996
$ws.Range("H134").Value = 5423.7666
$ws.Range("I134").Value = 5850.591
$ws.Range("K134").Value = 17551.773
$ws.Range("M134").Value = -15016.773

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2465.5334
$ws.Range("I31").Value = 995.28125
$ws.Range("J31").Value = 6084.615
$ws.Range("K31").Value = 995.28125
$ws.Range("L31").Value = 6084.615
$ws.Range("M31").Value = -700.28125
$ws.Range("N31").Value = -6674.615
$ws.Range("H34").Value = 2465.5334
$ws.Range("I34").Value = 995.28125
$ws.Range("J34").Value = 6084.615
$ws.Range("K34").Value = 995.28125
$ws.Range("L34").Value = 6084.615
$ws.Range("M34").Value = -793.28125
$ws.Range("N34").Value = -6488.615
$ws.Range("H58").Value = 3420.392
$ws.Range("I58").Value = 4411.5557
$ws.Range("J58").Value = 2305.3333
$ws.Range("K58").Value = 4411.5557
$ws.Range("L58").Value = 2305.3333
$ws.Range("M58").Value = -4208.5557
$ws.Range("N58").Value = -2711.3333
$ws.Range("H94").Value = 4979.591
$ws.Range("I94").Value = 1268
$ws.Range("J94").Value = 8691.182000000001
$ws.Range("K94").Value = 1268
$ws.Range("L94").Value = 8691.182000000001
$ws.Range("M94").Value = -817
$ws.Range("N94").Value = -9593.182000000001
$ws.Range("H99").Value = 336137.34
$ws.Range("I99").Value = 1000012
$ws.Range("J99").Value = 4200
$ws.Range("K99").Value = 1000012
$ws.Range("L99").Value = 4200
$ws.Range("M99").Value = -998514
$ws.Range("N99").Value = -7196
$ws.Range("H126").Value = 336137.34
$ws.Range("I126").Value = 1000012
$ws.Range("J126").Value = 4200
$ws.Range("K126").Value = 3000036
$ws.Range("L126").Value = 12600
$ws.Range("M126").Value = -2997566
$ws.Range("N126").Value = -17540
$ws.Range("H132").Value = 1703.4318
$ws.Range("I132").Value = 853.63336
$ws.Range("J132").Value = 3524.4285
$ws.Range("K132").Value = 2560.90008
$ws.Range("L132").Value = 10573.2855
$ws.Range("M132").Value = -30.90008000000034
$ws.Range("N132").Value = -15633.2855
$ws.Range("H134").Value = 1729.2941
$ws.Range("I134").Value = 1047.0588
$ws.Range("J134").Value = 2411.5293
$ws.Range("K134").Value = 3141.1764
$ws.Range("L134").Value = 7234.5879
$ws.Range("M134").Value = -606.1764000000003
$ws.Range("N134").Value = -12304.5879
$ws.Range("H136").Value = 3420.392
$ws.Range("I136").Value = 4411.5557
$ws.Range("J136").Value = 2305.3333
$ws.Range("K136").Value = 13234.6671
$ws.Range("L136").Value = 6915.999899999999
$ws.Range("M136").Value = -10684.6671
$ws.Range("N136").Value = -12015.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 8400
$ws.Range("I10").Value = 8400
$ws.Range("K10").Value = 8400
$ws.Range("M10").Value = -8231
$ws.Range("H102").Value = 21000
$ws.Range("I102").Value = 20000
$ws.Range("J102").Value = 21500
$ws.Range("K102").Value = 20000
$ws.Range("L102").Value = 21500
$ws.Range("M102").Value = -18378
$ws.Range("N102").Value = -24744
$ws.Range("H122").Value = 1581.75
$ws.Range("I122").Value = 1121.3
$ws.Range("J122").Value = 2349.1667
$ws.Range("K122").Value = 3363.9
$ws.Range("L122").Value = 7047.500100000001
$ws.Range("M122").Value = -913.8999999999996
$ws.Range("N122").Value = -11947.5001
$ws.Range("H126").Value = 2024.3529
$ws.Range("I126").Value = 1867.2
$ws.Range("J126").Value = 2248.8572
$ws.Range("K126").Value = 5601.6
$ws.Range("L126").Value = 6746.571599999999
$ws.Range("M126").Value = -3131.6
$ws.Range("N126").Value = -11686.5716
$ws.Range("H132").Value = 3722.3157
$ws.Range("I132").Value = 3113.923
$ws.Range("J132").Value = 5040.5
$ws.Range("K132").Value = 9341.769
$ws.Range("L132").Value = 15121.5
$ws.Range("M132").Value = -6811.769
$ws.Range("N132").Value = -20181.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2597.0435
$ws.Range("I122").Value = 2282.8
$ws.Range("J122").Value = 3186.25
$ws.Range("K122").Value = 6848.400000000001
$ws.Range("L122").Value = 9558.75
$ws.Range("M122").Value = -4398.400000000001
$ws.Range("N122").Value = -14458.75
$ws.Range("H132").Value = 11635052
$ws.Range("I132").Value = 20835930
$ws.Range("J132").Value = 12891.947
$ws.Range("K132").Value = 62507790
$ws.Range("L132").Value = 38675.841
$ws.Range("M132").Value = -62505260
$ws.Range("N132").Value = -43735.841
$ws.Range("H136").Value = 3787.8909
$ws.Range("I136").Value = 1882.1143
$ws.Range("J136").Value = 7123
$ws.Range("K136").Value = 5646.3429
$ws.Range("L136").Value = 21369
$ws.Range("M136").Value = -3096.3429
$ws.Range("N136").Value = -26469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3502.5
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 3005
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 9015
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -13915
$ws.Range("H123").Value = 18639.5
$ws.Range("J123").Value = 18639.5
$ws.Range("L123").Value = 18639.5
$ws.Range("N123").Value = -28439.5
$ws.Range("H126").Value = 1042.7727
$ws.Range("I126").Value = 662.2143
$ws.Range("J126").Value = 1708.75
$ws.Range("K126").Value = 1986.6429
$ws.Range("L126").Value = 5126.25
$ws.Range("M126").Value = 483.3571000000002
$ws.Range("N126").Value = -10066.25
$ws.Range("H132").Value = 24392278
$ws.Range("I132").Value = 31251386
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 93754158
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -93751628
$ws.Range("N132").Value = -18059
$ws.Range("H136").Value = 12183861
$ws.Range("I136").Value = 21301466
$ws.Range("K136").Value = 63904398
$ws.Range("M136").Value = -63901848
